# Updated solution for Tutorial 6
# Change the date format in column A from DD/MM/YYYY to DD-MM-YYYY
# and update the attendance counters (D, E, G, H) for a few rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New dates (text, dash-separated) for rows 3..21
$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

# Rows whose new date text is day/month-ambiguous (day <= 12) need to be
# pinned to Text format first, otherwise Excel auto-converts the typed
# string into a date serial number instead of keeping literal text.
$ambiguousRows = @(4, 5, 6, 7, 13, 14, 15, 16)

foreach ($r in $dates.Keys) {
    $cell = $ws.Cells.Item($r, 1)
    if ($ambiguousRows -contains $r) {
        $cell.NumberFormat = "@"
        $cell.Value = $dates[$r]
        # Restore the default ("Normal") cell style so no lingering
        # number-format override is left behind on the cell.
        $cell.Style = "Normal"
    } else {
        $cell.Value = $dates[$r]
    }
}

# Updated attendance counts: Total(D), Real(E), Duplicate(F), Invalid(G), Absent(H)
$ws.Cells.Item(3, 4).Value = 1   # D3
$ws.Cells.Item(3, 7).Value = 1   # G3

$ws.Cells.Item(4, 4).Value = 1   # D4
$ws.Cells.Item(4, 5).Value = 1   # E4
$ws.Cells.Item(4, 8).Value = 0   # H4

$ws.Cells.Item(5, 4).Value = 1   # D5
$ws.Cells.Item(5, 5).Value = 1   # E5
$ws.Cells.Item(5, 8).Value = 0   # H5

$ws.Cells.Item(10, 4).Value = 1  # D10
$ws.Cells.Item(10, 5).Value = 1  # E10
$ws.Cells.Item(10, 8).Value = 0  # H10
